# ---------------------------------------------------------------------------
# Applies the commit:
#   1) Merge the two split runs "SimpleKMeans = " + "4"/"6" into single runs
#      inside the first summary table ("Clustering instances based on
#      “is_rising” class" table is Tables(3) -- the one with the
#      Cobweb/EM/SimpleKMeans header row).
#   2) Append two new "Numerical prediction" sections (a label paragraph
#      pair + a 2-column/7-row results table) at the end of the document,
#      right after the existing last table and before the document's
#      trailing empty paragraph / section break.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: merge "SimpleKMeans = " / "4" and "SimpleKMeans = " / "6" runs
# ---------------------------------------------------------------------------

$kmTable = $d.Tables(3)

$cell5 = $kmTable.Cell(1, 5)
$cell5.Range.Find.Execute("SimpleKMeans = 4", $true, $false, $false, $false, `
    $false, $true, 1, $false, "SimpleKMeans = 4", 2) | Out-Null

$cell6 = $kmTable.Cell(1, 6)
$cell6.Range.Find.Execute("SimpleKMeans = 6", $true, $false, $false, $false, `
    $false, $true, 1, $false, "SimpleKMeans = 6", 2) | Out-Null

# ---------------------------------------------------------------------------
# Part 2: append the two new "Numerical prediction" sections
# ---------------------------------------------------------------------------

function Add-NumericPredictionSection($d, $className, $rows, $rowHeights, $tableWidthDxa, $colWidthDxa) {
    # Blank separator paragraph, mirrors "<w:p/>" in the target markup.
    $trailing = $d.Paragraphs($d.Paragraphs.Count)
    $trailing.Range.InsertParagraphBefore()

    # Label paragraph: Numerical prediction based on "<class>" class:
    $trailing = $d.Paragraphs($d.Paragraphs.Count)
    $trailing.Range.InsertParagraphBefore()
    $labelPara = $d.Paragraphs($d.Paragraphs.Count - 1)
    $labelPara.Range.Text = "Numerical prediction based on " + [char]8220 + $className + [char]8221 + " class:"

    # Results table, inserted right before the trailing empty paragraph.
    $trailing = $d.Paragraphs($d.Paragraphs.Count)
    $insertRange = $d.Range($trailing.Range.Start, $trailing.Range.Start)
    $newTable = $d.Tables.Add($insertRange, $rows.Count, 2)

    $newTable.Style = "Table Grid"
    $newTable.ApplyStyleHeadingRows = $true
    $newTable.ApplyStyleLastRow = $false
    $newTable.ApplyStyleFirstColumn = $true
    $newTable.ApplyStyleLastColumn = $false
    $newTable.ApplyStyleRowBands = $true
    $newTable.ApplyStyleColumnBands = $false

    if ($tableWidthDxa -ne $null) {
        $newTable.PreferredWidthType = 3
        $newTable.PreferredWidth = $tableWidthDxa / 20.0
    }

    $newTable.Columns(1).Width = $colWidthDxa / 20.0
    $newTable.Columns(2).Width = $colWidthDxa / 20.0

    for ($r = 1; $r -le $rows.Count; $r++) {
        $left = $rows[$r - 1][0]
        $right = $rows[$r - 1][1]

        if ($left -ne "") {
            $newTable.Cell($r, 1).Range.Text = $left
        }
        $newTable.Cell($r, 2).Range.Text = $right

        if ($rowHeights -ne $null) {
            $newTable.Rows($r).Height = $rowHeights[$r - 1] / 20.0
        }
    }

    return $newTable
}

$mmYearRows = @(
    @("", "Linear Regression"),
    @("Cross-Validation: 5", "58%"),
    @("Cross-Validation: 10", "59%"),
    @("Cross-Validation: 20", "60%"),
    @("Percent Split: 50%", "51%"),
    @("Percent Split: 66%", "48%"),
    @("Percent Split: 75%", "50%")
)
$mmYearHeights = @(259, 259, 271, 259, 259, 259, 259)

Add-NumericPredictionSection $d "msl_trend_mm_year" $mmYearRows $mmYearHeights 9596 4798 | Out-Null

$ftCenturyRows = @(
    @("", "Linear Regression"),
    @("Cross-Validation: 5", "58%"),
    @("Cross-Validation: 10", "59%"),
    @("Cross-Validation: 20", "60%"),
    @("Percent Split: 50%", "51%"),
    @("Percent Split: 66%", "50%"),
    @("Percent Split: 75%", "50%")
)

Add-NumericPredictionSection $d "msl_trend_ft_century" $ftCenturyRows $null $null 4788 | Out-Null

Write-Host "Edit complete."
